# Add "clutch" button support:
#  - Typography sheet: new bigger Consolas variant (Typography_01, size 150)
#    used for the clutch gauge readout, plus bump Default font size 20 -> 25.
#  - Translation sheet: re-home several existing text IDs onto the new
#    Typography_01 font, drop two now-unused rows, and append the new text
#    IDs for the clutch button and the various unit/label strings
#    (rpm, km/h, bar, degC, gear, ...).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Typography sheet
# ---------------------------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")

# Default font size 20 -> 25
$typo.Range("D4").Value2 = 25

# New row: Typography_01 / consola.ttf / 150 / 4 / ? / 32-127 / 32-127
$typo.Range("B8").Value2 = "Typography_01"
$typo.Range("C8").Value2 = "consola.ttf"
$typo.Range("D8").Value2 = 150
$typo.Range("E8").Value2 = 4
$typo.Range("F8").Value2 = "?"
$typo.Range("G8").Value2 = "32-127"
$typo.Range("H8").Value2 = "32-127"

# ---------------------------------------------------------------------
# Translation sheet
# ---------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Translation")

# Rows 4-7: switch typography name from Typography_00 to Typography_01
$tr.Range("C4").Value2 = "Typography_01"
$tr.Range("C5").Value2 = "Typography_01"
$tr.Range("C6").Value2 = "Typography_01"
$tr.Range("C7").Value2 = "Typography_01"

# Old rows 8 (SingleUseId9) and 9 (SingleUseId10) are dropped; rows 10-15
# shift up to become the new rows 8-13 (values unchanged other than the
# row position).
$tr.Range("B8").Value2 = "SingleUseId11"

$tr.Range("B9").Value2 = "SingleUseId13"
$tr.Range("D9").Value2 = "Right"
$tr.Range("E9").Value2 = "<value>"

$tr.Range("B10").Value2 = "SingleUseId14"
$tr.Range("D10").Value2 = "Left"
$tr.Range("E10").NumberFormat = "@"
$tr.Range("E10").Value2 = "3"

$tr.Range("B11").Value2 = "SingleUseId15"

$tr.Range("B12").Value2 = "SingleUseId16"

$tr.Range("B13").Value2 = "SingleUseId20"
$tr.Range("D13").Value2 = "Left"
$tr.Range("E13").NumberFormat = "@"
$tr.Range("E13").Value2 = "99.9"

# New row 14: reuses the old row-14 layout but on font Typography_01
$tr.Range("B14").Value2 = "SingleUseId21"
$tr.Range("C14").Value2 = "Typography_01"
$tr.Range("D14").Value2 = "Right"
$tr.Range("E14").Value2 = "<value>"

# New row 15
$tr.Range("B15").Value2 = "SingleUseId22"
$tr.Range("C15").Value2 = "Typography_01"
$tr.Range("D15").Value2 = "Left"
$tr.Range("E15").NumberFormat = "@"
$tr.Range("E15").Value2 = "999"

# New row 16 (brand new, same shape as old row 4/6/8/10/11)
$tr.Range("B16").Value2 = "SingleUseId23"
$tr.Range("C16").Value2 = "Typography_00"
$tr.Range("D16").Value2 = "Right"
$tr.Range("E16").Value2 = "<value>"
$tr.Range("F16").Value2 = "LTR"

# New row 17
$tr.Range("B17").Value2 = "SingleUseId24"
$tr.Range("C17").Value2 = "Typography_00"
$tr.Range("D17").Value2 = "Left"
$tr.Range("E17").NumberFormat = "@"
$tr.Range("E17").Value2 = "99.9"
$tr.Range("F17").Value2 = "LTR"

# New row 18: clutch button text - "rpm"
$tr.Range("B18").Value2 = "SingleUseId25"
$tr.Range("C18").Value2 = "Default"
$tr.Range("D18").Value2 = "Left"
$tr.Range("E18").Value2 = "rpm"
$tr.Range("F18").Value2 = "LTR"

# New row 19: "km/h"
$tr.Range("B19").Value2 = "SingleUseId26"
$tr.Range("C19").Value2 = "Default"
$tr.Range("D19").Value2 = "Left"
$tr.Range("E19").Value2 = "km/h"
$tr.Range("F19").Value2 = "LTR"

# New row 20: "bar"
$tr.Range("B20").Value2 = "SingleUseId27"
$tr.Range("C20").Value2 = "Default"
$tr.Range("D20").Value2 = "Left"
$tr.Range("E20").Value2 = "bar"
$tr.Range("F20").Value2 = "LTR"

# New row 21: "bar"
$tr.Range("B21").Value2 = "SingleUseId28"
$tr.Range("C21").Value2 = "Default"
$tr.Range("D21").Value2 = "Left"
$tr.Range("E21").Value2 = "bar"
$tr.Range("F21").Value2 = "LTR"

# New row 22: degrees C
$tr.Range("B22").Value2 = "SingleUseId29"
$tr.Range("C22").Value2 = "Default"
$tr.Range("D22").Value2 = "Left"
$tr.Range("E22").Value2 = [char]0x00B0 + "C"
$tr.Range("F22").Value2 = "LTR"

# New row 23: degrees C
$tr.Range("B23").Value2 = "SingleUseId30"
$tr.Range("C23").Value2 = "Default"
$tr.Range("D23").Value2 = "Left"
$tr.Range("E23").Value2 = [char]0x00B0 + "C"
$tr.Range("F23").Value2 = "LTR"

# New row 24: "gear"
$tr.Range("B24").Value2 = "SingleUseId31"
$tr.Range("C24").Value2 = "Default"
$tr.Range("D24").Value2 = "Left"
$tr.Range("E24").Value2 = "gear"
$tr.Range("F24").Value2 = "LTR"
